# Horarios actualizados Linea 141 - 313
# New scrape timestamp 06:15:33 merged into the three schedule sheets:
# new rows are inserted (keeping the existing rows sorted by Hora_Llegada)
# and the header "Ultima actualizacion" / "Total filas" cells are refreshed.

$wb = $excel.ActiveWorkbook

function Set-Row($ws, $row, $horaScrap, $horaLlegada, $linea, $minutos, $parada) {
    $ws.Cells.Item($row, 1).Value = $horaScrap
    $ws.Cells.Item($row, 2).Value = $horaLlegada
    $ws.Cells.Item($row, 3).Value = $linea
    $ws.Cells.Item($row, 4).Value = $minutos
    $ws.Cells.Item($row, 5).Value = $parada
}

# ---------------------------------------------------------------
# Sheet 1: LP1912  (40 -> 45 data rows, A1:E45 -> A1:E50)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 06:15:33"
$ws1.Cells.Item(3,1).Value = "Total filas: 45"

$ws1.Rows.Item(27).Insert()
Set-Row $ws1 27 "06:15:33" "06:15" "225_HARAS DEL SUR" 0 "LP1912"

$ws1.Rows.Item(35).Insert()
Set-Row $ws1 35 "06:15:33" "07:01" "16_SANTA ANA" 46 "LP1912"

$ws1.Rows.Item(42).Insert()
Set-Row $ws1 42 "06:15:33" "07:23" "10_OLMOS" 68 "LP1912"

Set-Row $ws1 49 "06:15:33" "08:07" "16_SANTA ANA" 112 "LP1912"
Set-Row $ws1 50 "06:15:33" "08:12" "15_ABASTO" 117 "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215  (only the "Última actualización" stamp changes)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 06:15:33"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173  (9 -> 11 data rows, A1:E14 -> A1:E16)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 06:15:33"
$ws3.Cells.Item(3,1).Value = "Total filas: 11"

$ws3.Rows.Item(14).Insert()
Set-Row $ws3 14 "06:15:33" "07:00" "215B_LP-P MOR-1 Y 57" 45 "L6173"

Set-Row $ws3 16 "06:15:33" "08:06" "215C_LA PLATA" 111 "L6203"
